# Update "展览" (Exhibitions) and "全部类型" (All Types) sheets with refreshed
# data snapshot values (增量更新 "想去人数" and sold-out status).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new F column (想去人数) value
$fUpdates = @{
    3  = 95
    4  = 391
    5  = 11555
    6  = 782
    7  = 113
    8  = 14
    14 = 50
    16 = 34
    17 = 328
    18 = 1322
    19 = 73
    20 = 899
}

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # G2: 已售罄 -> 不可售
    $ws.Range("G2").Value = "不可售"

    foreach ($row in $fUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $fUpdates[$row]
    }
}
